# Update the "dSF" (column F) values for the rows that were repulled/recalculated.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    5  = -1
    11 = 0
    13 = -2
    14 = -1
    15 = -4
    21 = 0
    30 = -6
    32 = -2
    33 = -2
    37 = -2
    42 = 1
    43 = -1
    44 = -3
    45 = -2
    48 = -1
    54 = -1
    57 = -1
    60 = -4
    63 = -3
    67 = 6
    69 = -8
    70 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
